$d = $word.ActiveDocument

# The new bullet goes after the last paragraph in the document
# ("Assisted with the data formatting in the display"), inheriting its
# ListParagraph style / numbering (numId 2, ilvl 0).
$count = $d.Paragraphs.Count
$last = $d.Paragraphs($count)
$last.Range.InsertParagraphAfter()

$newCount = $d.Paragraphs.Count
$newPara = $d.Paragraphs($newCount)
$newPara.Range.Text = "Tested the program" + [char]0x2019 + "s functionality"
